$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Row 6: add Lucas / 80% ---
$ws.Range("B6").Value2 = "Lucas"
$ws.Range("C6").Value2 = 0.8
$ws.Range("C6").NumberFormat = "0%"

# --- Row 20: add 100% ---
$ws.Range("C20").Value2 = 1
$ws.Range("C20").NumberFormat = "0%"

# --- Row 25: add Agustina ---
$ws.Range("B25").Value2 = "Agustina"

# --- Row 35: add Agustina (trailing space) / 100%, then two new rows 36/37 ---
# Copy row 35 (which already carries the s="1" style on column A)
# down into two freshly inserted rows so the style index is preserved
# exactly as Excel would when duplicating a formatted row.
$ws.Rows("35").Copy() | Out-Null
$ws.Rows("36:37").Insert(-4121) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B35").Value2 = "Agustina "
$ws.Range("C35").Value2 = 1
$ws.Range("C35").NumberFormat = "0%"

$ws.Range("A36").Value2 = "producto comodin"

$ws.Range("A37").Value2 = "cambiar precio de producto"
$ws.Range("A37").ClearFormats() | Out-Null

# --- Update selection to reflect the new active cell ---
$ws.Range("C36").Select() | Out-Null
